$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("A58").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-09 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.01602883595714602
$ws.Range("E2").Value = 0.003042921204356341
$ws.Range("D3").Value = 0.05228124033483193
$ws.Range("E3").Value = 0.02209559603552269
$ws.Range("D4").Value = 0.01494324713137732
$ws.Range("E4").Value = 0.01563063349663851
$ws.Range("D5").Value = 0.00952290916000988
$ws.Range("E5").Value = 0.0006656017039403661
$ws.Range("D6").Value = 0.01535924056067142
$ws.Range("E6").Value = 0.0002104672372669203
$ws.Range("D7").Value = 0.02044541468025789
$ws.Range("E7").Value = -0.006484044878332984
$ws.Range("D8").Value = 0.004177808763009659
$ws.Range("E8").Value = -0.02023152081563306
$ws.Range("D9").Value = 0.006389172279878832
$ws.Range("E9").Value = -0.002321428571428585
$ws.Range("D10").Value = 0.01389388263058829
$ws.Range("E10").Value = 0.00480384307445969
$ws.Range("D11").Value = 0.00871855800691799
$ws.Range("E11").Value = -0.00625227190112676
$ws.Range("D12").Value = 0.01450494230682016
$ws.Range("E12").Value = 0.007865757734661472
$ws.Range("D13").Value = 0.003167964588773254
$ws.Range("E13").Value = 0.02591036414565839
$ws.Range("D14").Value = 0.006114716764284831
$ws.Range("E14").Value = -0.0114025085518813
$ws.Range("D15").Value = 0.01425672803451852
$ws.Range("E15").Value = 0.007478081485301491
$ws.Range("D16").Value = 0.01031537569204545
$ws.Range("E16").Value = 0.001244296972210845
$ws.Range("D17").Value = 0.02130532247524774
$ws.Range("E17").Value = -0.004417966396679796
$ws.Range("D18").Value = 0.008240701163485385
$ws.Range("E18").Value = 0.01130673558391226
$ws.Range("D19").Value = 0.01652767311828339
$ws.Range("E19").Value = -0.01055408970976257
$ws.Range("D20").Value = 0.01157467013917335
$ws.Range("E20").Value = 0.01029516455834845
$ws.Range("D21").Value = 0.007035695665342762
$ws.Range("E21").Value = 0.001333333333333409
$ws.Range("D22").Value = 0.01309602840367386
$ws.Range("E22").Value = 0.003794552107331484
$ws.Range("D23").Value = 0.01867977199148503
$ws.Range("E23").Value = 0.01983508932661482
$ws.Range("D24").Value = 0.009466496825395869
$ws.Range("E24").Value = 0.01048543689320391
$ws.Range("D25").Value = 0.02111913007869869
$ws.Range("E25").Value = -0.0003901677721420693
$ws.Range("D26").Value = 0.01151515195692329
$ws.Range("E26").Value = 0.0168766616576117
$ws.Range("D27").Value = 0.0232875188058282
$ws.Range("E27").Value = -0.008982035928143728
$ws.Range("D28").Value = 0.05783975683368088
$ws.Range("E28").Value = 0.02021325559987708
$ws.Range("D29").Value = 0.02132915510200602
$ws.Range("E29").Value = 0.004038590980480095
$ws.Range("D30").Value = 0.03253914168252577
$ws.Range("E30").Value = -0.002337540906966029
$ws.Range("D31").Value = 0.01645649216123672
$ws.Range("E31").Value = 0.001638877012968498
$ws.Range("D32").Value = 0.01361350065062754
$ws.Range("E32").Value = 0.001238499646142976
$ws.Range("D33").Value = 0.02174372237677336
$ws.Range("E33").Value = 0.0001049428061705715
$ws.Range("D34").Value = 0.04279281242204678
$ws.Range("E34").Value = 0.008993836733424354
$ws.Range("D35").Value = 0.01092244213560238
$ws.Range("E35").Value = -0.002785515320334109
$ws.Range("D36").Value = 0.009415155262432556
$ws.Range("E36").Value = 0.01357210179076374
$ws.Range("D37").Value = 0.01183945949630821
$ws.Range("E37").Value = -0.01084119654687798
$ws.Range("D38").Value = 0.007257225001832609
$ws.Range("E38").Value = 0.003013232018865386
$ws.Range("D39").Value = 0.01170112250721262
$ws.Range("E39").Value = 0.007304785894206578
$ws.Range("D40").Value = 0.01743857386068556
$ws.Range("E40").Value = 0.003467539973030176
$ws.Range("D41").Value = 0.01684827265589873
$ws.Range("E41").Value = 0.00393137955682632
$ws.Range("D42").Value = 0.03436652101614587
$ws.Range("E42").Value = 0.006363083052066676
$ws.Range("D43").Value = 0.0111711634850691
$ws.Range("E43").Value = 0.004204398447606783
$ws.Range("D44").Value = 0.02168293650161062
$ws.Range("E44").Value = 0.008246488445854183
$ws.Range("D45").Value = 0.01384738999301596
$ws.Range("E45").Value = -0.000803328664414682
$ws.Range("D46").Value = 0.008113107871824711
$ws.Range("E46").Value = 0.007523555055547693
$ws.Range("D47").Value = 0.01314800381309351
$ws.Range("E47").Value = 0.001171468240194429
$ws.Range("D48").Value = 0.009415186954755374
$ws.Range("E48").Value = 0.0199945469417433
$ws.Range("D49").Value = 0.01522372418830653
$ws.Range("E49").Value = 0.009680239820134995
$ws.Range("D50").Value = 0.008219467307198089
$ws.Range("E50").Value = 0.009438909281594077
$ws.Range("D51").Value = 0.01101580771862085
$ws.Range("E51").Value = 0.02689406362741886
$ws.Range("D52").Value = 0.008466825886783677
$ws.Range("E52").Value = 0.02083419113105767
$ws.Range("D53").Value = 0.1376235315243377
$ws.Range("E53").Value = -0.0001971220185293943
$ws.Range("D54").Value = 0.04400130407569926
$ws.Range("E54").Value = 0.007224194930819161
$ws.Range("E55").Value = 0.005422493049319188
